$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# New column H: header "Error_user" (row 1) and value "error_user" (row 2)
$ws.Range("H1").Value = "Error_user"
$ws.Range("H2").Value = "error_user"

# Match the direct formatting used by the other header cells (bold font,
# yellow fill, thin black border) so H1 reuses the existing header style.
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").Interior.Color = 65535
$ws.Range("H1").Borders.LineStyle = 1
$ws.Range("H1").Borders.Color = 0

# Column H width (matches the auto-fit width Excel computed for the new column)
$ws.Columns.Item(8).ColumnWidth = 8.75

# Move the active selection to F1
$ws.Range("F1").Select()
